$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Il1a"
$ws.Cells.Item(2,3).Value = "Il1r1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.038687
$ws.Cells.Item(2,8).Value = 0.116061
$ws.Cells.Item(2,9).Value = 0.006675966692470672
$ws.Cells.Item(2,10).Value = 0.006675966692470672
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 30.60409033333334
$ws.Cells.Item(2,14).Value = 91.81227100000001
$ws.Cells.Item(2,15).Value = 0.2776440509987301
$ws.Cells.Item(2,16).Value = 0.27764405099873
$ws.Cells.Item(2,17).Value = 1.183980442725667
$ws.Cells.Item(2,18).Value = 10.655823984531
$ws.Cells.Item(2,19).Value = 0.001853542436830151
$ws.Cells.Item(2,20).Value = 0.00185354243683015

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Il1a"
$ws.Cells.Item(3,3).Value = "Il1r1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.038687
$ws.Cells.Item(3,8).Value = 0.116061
$ws.Cells.Item(3,9).Value = 0.006675966692470672
$ws.Cells.Item(3,10).Value = 0.006675966692470672
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 72.39518466666667
$ws.Cells.Item(3,14).Value = 217.185554
$ws.Cells.Item(3,15).Value = 0.6567779706806668
$ws.Cells.Item(3,16).Value = 0.6567779706806667
$ws.Cells.Item(3,17).Value = 2.800752509199333
$ws.Cells.Item(3,18).Value = 25.206772582794
$ws.Cells.Item(3,19).Value = 0.004384627856612611
$ws.Cells.Item(3,20).Value = 0.00438462785661261

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Il1a"
$ws.Cells.Item(4,3).Value = "Il1r1"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.038687
$ws.Cells.Item(4,8).Value = 0.116061
$ws.Cells.Item(4,9).Value = 0.006675966692470672
$ws.Cells.Item(4,10).Value = 0.006675966692470672
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.371039
$ws.Cells.Item(4,14).Value = 1.113117
$ws.Cells.Item(4,15).Value = 0.003366111193519582
$ws.Cells.Item(4,16).Value = 0.003366111193519582
$ws.Cells.Item(4,17).Value = 0.014354385793
$ws.Cells.Item(4,18).Value = 0.129189472137
$ws.Cells.Item(4,19).Value = [double]"2.247204621108943E-05"
$ws.Cells.Item(4,20).Value = [double]"2.247204621108943E-05"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Il1a"
$ws.Cells.Item(5,3).Value = "Il1r1"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.038687
$ws.Cells.Item(5,8).Value = 0.116061
$ws.Cells.Item(5,9).Value = 0.006675966692470672
$ws.Cells.Item(5,10).Value = 0.006675966692470672
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 6.687497333333333
$ws.Cells.Item(5,14).Value = 20.062492
$ws.Cells.Item(5,15).Value = 0.06066979382319835
$ws.Cells.Item(5,16).Value = 0.06066979382319834
$ws.Cells.Item(5,17).Value = 0.2587192093346666
$ws.Cells.Item(5,18).Value = 2.328472884012
$ws.Cells.Item(5,19).Value = 0.0004050295228027351
$ws.Cells.Item(5,20).Value = 0.000405029522802735

# Row 6
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Il1a"
$ws.Cells.Item(6,3).Value = "Il1r1"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.038687
$ws.Cells.Item(6,8).Value = 0.116061
$ws.Cells.Item(6,9).Value = 0.006675966692470672
$ws.Cells.Item(6,10).Value = 0.006675966692470672
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.1699793333333333
$ws.Cells.Item(6,14).Value = 0.509938
$ws.Cells.Item(6,15).Value = 0.001542073303885386
$ws.Cells.Item(6,16).Value = 0.001542073303885386
$ws.Cells.Item(6,17).Value = 0.006575990468666667
$ws.Cells.Item(6,18).Value = 0.059183914218
$ws.Cells.Item(6,19).Value = [double]"1.029483001408704E-05"
$ws.Cells.Item(6,20).Value = [double]"1.029483001408704E-05"

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Il1a"
$ws.Cells.Item(7,3).Value = "Il1r1"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.929373
$ws.Cells.Item(7,8).Value = 8.788119
$ws.Cells.Item(7,9).Value = 0.5055030521317986
$ws.Cells.Item(7,10).Value = 0.5055030521317986
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 30.60409033333334
$ws.Cells.Item(7,14).Value = 91.81227100000001
$ws.Cells.Item(7,15).Value = 0.2776440509987301
$ws.Cells.Item(7,16).Value = 0.27764405099873
$ws.Cells.Item(7,17).Value = 89.65079591202768
$ws.Cells.Item(7,18).Value = 806.8571632082491
$ws.Cells.Item(7,19).Value = 0.1403499151860948
$ws.Cells.Item(7,20).Value = 0.1403499151860948

# Row 8
$ws.Cells.Item(8,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,2).Value = "Il1a"
$ws.Cells.Item(8,3).Value = "Il1r1"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.929373
$ws.Cells.Item(8,8).Value = 8.788119
$ws.Cells.Item(8,9).Value = 0.5055030521317986
$ws.Cells.Item(8,10).Value = 0.5055030521317986
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 72.39518466666667
$ws.Cells.Item(8,14).Value = 217.185554
$ws.Cells.Item(8,15).Value = 0.6567779706806668
$ws.Cells.Item(8,16).Value = 0.6567779706806667
$ws.Cells.Item(8,17).Value = 212.0724992925473
$ws.Cells.Item(8,18).Value = 1908.652493632926
$ws.Cells.Item(8,19).Value = 0.332003268752006
$ws.Cells.Item(8,20).Value = 0.3320032687520059

# Row 9
$ws.Cells.Item(9,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,2).Value = "Il1a"
$ws.Cells.Item(9,3).Value = "Il1r1"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.929373
$ws.Cells.Item(9,8).Value = 8.788119
$ws.Cells.Item(9,9).Value = 0.5055030521317986
$ws.Cells.Item(9,10).Value = 0.5055030521317986
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.371039
$ws.Cells.Item(9,14).Value = 1.113117
$ws.Cells.Item(9,15).Value = 0.003366111193519582
$ws.Cells.Item(9,16).Value = 0.003366111193519582
$ws.Cells.Item(9,17).Value = 1.086911628547
$ws.Cells.Item(9,18).Value = 9.782204656923
$ws.Cells.Item(9,19).Value = 0.00170157948213916
$ws.Cells.Item(9,20).Value = 0.00170157948213916

# Row 10
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Il1a"
$ws.Cells.Item(10,3).Value = "Il1r1"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.929373
$ws.Cells.Item(10,8).Value = 8.788119
$ws.Cells.Item(10,9).Value = 0.5055030521317986
$ws.Cells.Item(10,10).Value = 0.5055030521317986
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 6.687497333333333
$ws.Cells.Item(10,14).Value = 20.062492
$ws.Cells.Item(10,15).Value = 0.06066979382319835
$ws.Cells.Item(10,16).Value = 0.06066979382319834
$ws.Cells.Item(10,17).Value = 19.59017412583867
$ws.Cells.Item(10,18).Value = 176.311567132548
$ws.Cells.Item(10,19).Value = 0.03066876594983371
$ws.Cells.Item(10,20).Value = 0.0306687659498337

# Row 11
$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Il1a"
$ws.Cells.Item(11,3).Value = "Il1r1"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.929373
$ws.Cells.Item(11,8).Value = 8.788119
$ws.Cells.Item(11,9).Value = 0.5055030521317986
$ws.Cells.Item(11,10).Value = 0.5055030521317986
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.1699793333333333
$ws.Cells.Item(11,14).Value = 0.509938
$ws.Cells.Item(11,15).Value = 0.001542073303885386
$ws.Cells.Item(11,16).Value = 0.001542073303885386
$ws.Cells.Item(11,17).Value = 0.4979328696246667
$ws.Cells.Item(11,18).Value = 4.481395826622
$ws.Cells.Item(11,19).Value = 0.0007795227617250291
$ws.Cells.Item(11,20).Value = 0.0007795227617250289

# Row 12
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Il1a"
$ws.Cells.Item(12,3).Value = "Il1r1"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 2.826906
$ws.Cells.Item(12,8).Value = 8.480718
$ws.Cells.Item(12,9).Value = 0.4878209811757308
$ws.Cells.Item(12,10).Value = 0.4878209811757308
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 30.60409033333334
$ws.Cells.Item(12,14).Value = 91.81227100000001
$ws.Cells.Item(12,15).Value = 0.2776440509987301
$ws.Cells.Item(12,16).Value = 0.27764405099873
$ws.Cells.Item(12,17).Value = 86.514886587842
$ws.Cells.Item(12,18).Value = 778.6339792905781
$ws.Cells.Item(12,19).Value = 0.1354405933758051
$ws.Cells.Item(12,20).Value = 0.1354405933758051

# Row 13
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Il1a"
$ws.Cells.Item(13,3).Value = "Il1r1"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 2.826906
$ws.Cells.Item(13,8).Value = 8.480718
$ws.Cells.Item(13,9).Value = 0.4878209811757308
$ws.Cells.Item(13,10).Value = 0.4878209811757308
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 72.39518466666667
$ws.Cells.Item(13,14).Value = 217.185554
$ws.Cells.Item(13,15).Value = 0.6567779706806668
$ws.Cells.Item(13,16).Value = 0.6567779706806667
$ws.Cells.Item(13,17).Value = 204.654381905308
$ws.Cells.Item(13,18).Value = 1841.889437147772
$ws.Cells.Item(13,19).Value = 0.3203900740720482
$ws.Cells.Item(13,20).Value = 0.3203900740720482

# Row 14
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Il1a"
$ws.Cells.Item(14,3).Value = "Il1r1"
$ws.Cells.Item(14,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 2.826906
$ws.Cells.Item(14,8).Value = 8.480718
$ws.Cells.Item(14,9).Value = 0.4878209811757308
$ws.Cells.Item(14,10).Value = 0.4878209811757308
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.371039
$ws.Cells.Item(14,14).Value = 1.113117
$ws.Cells.Item(14,15).Value = 0.003366111193519582
$ws.Cells.Item(14,16).Value = 0.003366111193519582
$ws.Cells.Item(14,17).Value = 1.048892375334
$ws.Cells.Item(14,18).Value = 9.440031378005999
$ws.Cells.Item(14,19).Value = 0.001642059665169333
$ws.Cells.Item(14,20).Value = 0.001642059665169333

# Row 15
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Il1a"
$ws.Cells.Item(15,3).Value = "Il1r1"
$ws.Cells.Item(15,4).Value = "MuSCs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 2.826906
$ws.Cells.Item(15,8).Value = 8.480718
$ws.Cells.Item(15,9).Value = 0.4878209811757308
$ws.Cells.Item(15,10).Value = 0.4878209811757308
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 6.687497333333333
$ws.Cells.Item(15,14).Value = 20.062492
$ws.Cells.Item(15,15).Value = 0.06066979382319835
$ws.Cells.Item(15,16).Value = 0.06066979382319834
$ws.Cells.Item(15,17).Value = 18.904926336584
$ws.Cells.Item(15,18).Value = 170.144337029256
$ws.Cells.Item(15,19).Value = 0.02959599835056191
$ws.Cells.Item(15,20).Value = 0.0295959983505619

# Row 16
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Il1a"
$ws.Cells.Item(16,3).Value = "Il1r1"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 2.826906
$ws.Cells.Item(16,8).Value = 8.480718
$ws.Cells.Item(16,9).Value = 0.4878209811757308
$ws.Cells.Item(16,10).Value = 0.4878209811757308
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.1699793333333333
$ws.Cells.Item(16,14).Value = 0.509938
$ws.Cells.Item(16,15).Value = 0.001542073303885386
$ws.Cells.Item(16,16).Value = 0.001542073303885386
$ws.Cells.Item(16,17).Value = 0.480515597276
$ws.Cells.Item(16,18).Value = 4.324640375484
$ws.Cells.Item(16,19).Value = 0.0007522557121462699
$ws.Cells.Item(16,20).Value = 0.0007522557121462696
